# Creado grafico de tipos de modelo
# Insert a new "MAE" column before the existing "Tipo" column, shifting
# "Tipo" from D to E, and populate the new MAE values. Also refresh the
# MSE value in B2 with the latest computed figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the old column D ("Tipo") to E
# and carries over the header styling (bold/border/centered) from the
# neighboring header cell.
$ws.Columns.Item(4).Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "MAE"

# Updated MSE value (tiny refinement from re-running the model).
$ws.Range("B2").Value = 0.3638702225807681

# New MAE metric value for the row.
$ws.Range("D2").Value = 0.4682092966157618
